$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 25)
$ws.Range("D2").Value = 0.3066410164048881
$ws.Range("E2").Value = 0.3066410164048881

# Row 3 (Control 44)
$ws.Range("D3").Value = [double]"2.67179446293341E-13"
$ws.Range("E3").Value = [double]"2.67179446293341E-13"

# Row 4 (Control 40)
$ws.Range("D4").Value = 0.5444571126839456
$ws.Range("E4").Value = 0.5444571126839456

# Row 5 (Control 41)
$ws.Range("C5").Value = $false
$ws.Range("D5").Value = 0.9878328960091652
$ws.Range("E5").Value = 0.9878328960091652

# Row 6 (Control 43)
$ws.Range("D6").Value = 0.3422016693524062
$ws.Range("E6").Value = 0.3422016693524062

# Row 7 (MDD 42)
$ws.Range("D7").Value = 0.09375203912300566
$ws.Range("E7").Value = 0.9062479608769943

# Row 9 (MDD 20)
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = 0.04163085272417932
$ws.Range("E9").Value = 0.9583691472758207

# Row 10 (MDD 51)
$ws.Range("D10").Value = 0.599051495102949
$ws.Range("E10").Value = 0.400948504897051

# Row 11 (MDD 40)
$ws.Range("D11").Value = 0.002086354461286365
$ws.Range("E11").Value = 0.9979136455387136
$ws.Range("F11").Value = 1.821110963821411
$ws.Range("G11").Value = 0.5
